# update file with RFLP plate
#
# 1. Rename "PocHistone RLFP 012" -> "DONE PocHistone RLFP 012" (it's done now).
# 2. The previously-active sheet "PocHistone RLFP 014" is no longer the active
#    tab; its selection moves to G19.
# 3. "PocHistone RLFP 016" becomes the new active sheet/tab (selection N17
#    stays as-is).

$wb = $excel.ActiveWorkbook

$doneSheet = $wb.Worksheets.Item(9)
$doneSheet.Name = "DONE PocHistone RLFP 012"

$prevActive = $wb.Worksheets.Item(11)
$prevActive.Range("G19").Select()

$newActive = $wb.Worksheets.Item(13)
$newActive.Activate()
